$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Todo" worksheet right after "Summary" (i.e. before the
#    sheet currently in position 2 = GlobalVars). It becomes the new active
#    tab, matching the commit: sheetId 20, new rId2.
# ---------------------------------------------------------------------------
$todo = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$todo.Name = "Todo"
$todo.Range("A1").Value = "The wooden tree root, lizard, pile of trash and green liquid are not part of 3Objects.amb therefore they must be added there. Maybe smaller or even different versions."
$todo.Range("A2").Value = "the teleports to same map in Ship's End which use map index 0 seem to teleport to map 0. I guess index 372 is necessary there."

# ---------------------------------------------------------------------------
# 2. Fix two typos ("met" -> "meet") in the descriptions of the first two
#    rows of the "Maps" sheet.
# ---------------------------------------------------------------------------
$maps = $wb.Worksheets.Item("Maps")
$maps.Range("D2").Value = "You can meet Karl the architect here who can renovate your house, create a wind gate or a cave"
$maps.Range("D3").Value = "You can meet Ferdinand who sells you a cat and a dog, there is also a merchant for pet stuff"

# ---------------------------------------------------------------------------
# 3. Append the first maps for the sea quests (new rows 5-9 on "Maps").
# ---------------------------------------------------------------------------
$maps.Range("A5").Value = 459
$maps.Range("B5").Value = "Deine Höhle - Obere Ebene / Your Cave - Upper Level"
$maps.Range("C5").Value = "3D"
$maps.Range("D5").Value = "Small cave, Karl built for your"

$maps.Range("A6").Value = 460
$maps.Range("B6").Value = "Deine Höhle - Untere Ebene / Your Cave - Lower Level"
$maps.Range("C6").Value = "3D"
$maps.Range("D6").Value = "Small cave, Karl built for your"

$maps.Range("A7").Value = 370
$maps.Range("B7").Value = "Auge des Strudels / Eye of the vortex"
$maps.Range("C7").Value = "2D"
$maps.Range("D7").Value = "Map inside the vortex"

$maps.Range("A8").Value = 371
$maps.Range("B8").Value = "Höhle der Meerjungfrau / Cave of the mermaid"
$maps.Range("C8").Value = "2D"
$maps.Range("D8").Value = "You can meet the Mermaid here"

$maps.Range("A9").Value = 372
$maps.Range("B9").Value = "Ship's end"
$maps.Range("C9").Value = "3D"
$maps.Range("D9").Value = "Abandoned village inside the vortex"

# Widen column B on "Maps" to fit the new, longer names (closest achievable
# approximation of the author's manual resize to ~49.14 chars).
$maps.Columns.Item(2).ColumnWidth = 48.25

# Move the selection/active cell the way the author's session ended up.
$maps.Range("A10").Select() | Out-Null
$todo.Activate() | Out-Null
